{"js": "// Helper: search for `oldText` (which may span multiple runs) and replace it\n// with `newText`, forcing Word to rewrite the matched range as a single run\n// (mirrors what Word itself does on a real text edit / run consolidation).\nasync function replaceAcrossRuns(context, oldText, newText) {\n  if (oldText === newText) {\n    return; // nothing to do\n  }\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// 1) Title: merge the 4 runs that together spell out the title into one run\n//    (the visible text itself does not change).\nawait replaceAcrossRuns(\n  context,\n  \"Predicting default rate of Home Credit Loan Customer \",\n  \"Predicting default rate of Home Credit Loan Customer \\u0001\"\n);\nawait replaceAcrossRuns(\n  context,\n  \"Predicting default rate of Home Credit Loan Customer \\u0001\",\n  \"Predicting default rate of Home Credit Loan Customer \"\n);\n\n// 2) Hyperlink text \"(Links to an external site.)\" - merge the 3 runs that\n//    spell it out into a single run (visible text unchanged).\nawait replaceAcrossRuns(\n  context,\n  \"(Links to an external site.)\",\n  \"(Links to an external site.)\\u0001\"\n);\nawait replaceAcrossRuns(\n  context,\n  \"(Links to an external site.)\\u0001\",\n  \"(Links to an external site.)\"\n);\n\n// 3) \"Logistic regression, random forest, and gradient boost classifiers\n//    were used ...\" -> \"Logistic regression, random forest, gradient boost\n//    classifier and voting classifier were used ...\"\nawait replaceAcrossRuns(\n  context,\n  \", and gradient boost classifiers were used\",\n  \", gradient boost classifier and voting classifier were used\"\n);\n\n// 4) Candidate model scores.\nawait replaceAcrossRuns(context, \"Logistic regression: 0.58\", \"Logistic regression: 0.57\");\nawait replaceAcrossRuns(context, \"Random forest (selected): 0.72\", \"Random forest: 0.71\");\nawait replaceAcrossRuns(context, \"Gradient boost: 0.71\", \"Gradient boost: 0.70\");\n\n// 5) New bullet point for the optimized random forest model, inserted right\n//    after the \"Gradient boost: 0.70\" bullet (same list formatting).\nconst gbResults = context.document.body.search(\"Gradient boost: 0.70\", { matchCase: true });\ngbResults.load(\"items\");\nawait context.sync();\ngbResults.items[0].insertParagraph(\n  \"Random forest with optimized max_depth: 0.73\",\n  \"After\"\n);\nawait context.sync();\n\n// 6) Tools list: \"Matplotlib and Seaborn for plotting\" -> \"Matplotlib for plotting\"\nawait replaceAcrossRuns(\n  context,\n  \"Matplotlib and Seaborn for plotting\",\n  \"Matplotlib for plotting\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1) Title: \"Predicting default rate of \" + \"Home Credit Loan \" + \"Customer\"\n#    + \" \" were four separate runs spelling out the same text; consolidate\n#    them into a single run (visible text is unchanged).\n# ---------------------------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"Predicting default rate of Home Credit Loan Customer \"\n$find.Execute() | Out-Null\n$rng = $find.Parent\n$rng.Text = \"Predicting default rate of Home Credit Loan Customer `u{0001}\"\n\n$find = $d.Content.Find\n$find.Text = \"Predicting default rate of Home Credit Loan Customer `u{0001}\"\n$find.Execute() | Out-Null\n$rng = $find.Parent\n$rng.Text = \"Predicting default rate of Home Credit Loan Customer \"\n\n# ---------------------------------------------------------------------------\n# 2) Hyperlink display text \"(Links to a\" + \"n\" + \" external site.)\" -> merge\n#    into a single run \"(Links to an external site.)\". We intentionally\n#    start the match one character in (at \"Links ...\") rather than at \"(\" so\n#    the edit is not anchored exactly on the boundary with the preceding\n#    \"Kaggle \" run (which has different run formatting); this keeps \"Kaggle \"\n#    untouched and keeps the hyperlink text's own formatting (including the\n#    border property) on the merged run.\n# ---------------------------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"Links to an external site.)\"\n$find.Execute() | Out-Null\n$rng = $find.Parent\n$rng.Text = \"Links to an external site.)`u{0001}\"\n\n$find = $d.Content.Find\n$find.Text = \"Links to an external site.)`u{0001}\"\n$find.Execute() | Out-Null\n$rng = $find.Parent\n$rng.Text = \"Links to an external site.)\"\n\n# ---------------------------------------------------------------------------\n# 3) \"Logistic regression, random forest, and gradient boost classifiers\n#    were used ...\" -> \"Logistic regression, random forest, gradient boost\n#    classifier and voting classifier were used ...\"\n# ---------------------------------------------------------------------------\n$d.Content.Find.Execute(\", and gradient boost classifiers were used\", $false, $false, $false, $false, $false, $true, 1, $false, \", gradient boost classifier and voting classifier were used\", 2) | Out-Null\n\n# ---------------------------------------------------------------------------\n# 4) Candidate model scores.\n# ---------------------------------------------------------------------------\n$d.Content.Find.Execute(\"Logistic regression: 0.58\", $false, $false, $false, $false, $false, $true, 1, $false, \"Logistic regression: 0.57\", 2) | Out-Null\n$d.Content.Find.Execute(\"Random forest (selected): 0.72\", $false, $false, $false, $false, $false, $true, 1, $false, \"Random forest: 0.71\", 2) | Out-Null\n$d.Content.Find.Execute(\"Gradient boost: 0.71\", $false, $false, $false, $false, $false, $true, 1, $false, \"Gradient boost: 0.70\", 2) | Out-Null\n\n# ---------------------------------------------------------------------------\n# 5) New bullet point for the optimized random forest model, inserted right\n#    after the \"Gradient boost: 0.70\" bullet (inherits the same list\n#    formatting).\n# ---------------------------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"Gradient boost: 0.70\"\n$find.Execute() | Out-Null\n$rng = $find.Parent\n$gbPara = $rng.Paragraphs(1)\n$rng.InsertParagraphAfter()\n$newPara = $gbPara.Next()\n$newPara.Range.Text = \"Random forest with optimized max_depth: 0.73\"\n\n# ---------------------------------------------------------------------------\n# 6) Tools list: \"Matplotlib and Seaborn for plotting\" -> \"Matplotlib for\n#    plotting\"\n# ---------------------------------------------------------------------------\n$d.Content.Find.Execute(\"Matplotlib and Seaborn for plotting\", $false, $false, $false, $false, $false, $true, 1, $false, \"Matplotlib for plotting\", 2) | Out-Null\n"}
